# Termo de abertura de projeto - content update
# (version bump 2.0 -> 3.0, date bump 04/08/2023 -> 12/04/2025,
#  schedule/financial figures updated)

$d = $word.ActiveDocument

# 1) Version bump: "2.0" -> "3.0"
#    Covers both "GENESYS PETSHOP 2.0" (title) and the VERSAO table cell.
$d.Content.Find.Execute("2.0", $false, $false, $false, $false, $false, $true, 1, $false, "3.0", 2)

# 2) Date bump: "04/08/2023" -> "12/04/2025"
#    Covers both the DATA table cell and the "Atualizado - ..." line.
$d.Content.Find.Execute("04/08/2023", $false, $false, $false, $false, $false, $true, 1, $false, "12/04/2025", 2)

# 3) Cronograma: "3 meses (tres meses)" -> "3/5 meses (tres meses e meio)"
$d.Content.Find.Execute("Cronograma – 3 meses (três meses)", $false, $false, $false, $false, $false, $true, 1, $false, "Cronograma – 3/5 meses (três meses e meio)", 2)

# 4) Financeiro: "R$ 11.235,00" -> "R$ 13.107,5"
$d.Content.Find.Execute("Financeiro – R$ 11.235,00", $false, $false, $false, $false, $false, $true, 1, $false, "Financeiro – R$ 13.107,5", 2)

# 5) Tidy up: the now-empty paragraph right after the "Atualizado - ..." line
#    carries a stale paragraph-mark run format; re-applying its own style
#    clears that redundant formatting hint (matches a plain "Standard" pPr).
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq ([char]13)) {
        $prev = $paras.Item($i - 1).Range.Text
        if ($prev -like "*Atualizado*") {
            $p.Style = $p.Style
        }
    }
}
